$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    # Temporarily force text number format so Excel does not auto-convert
    # numeric-looking strings (e.g. "570.76") into numbers, matching the
    # original cells which are stored as text (inlineStr). Restore the
    # original number format afterwards so no other formatting changes.
    $cell = $ws.Range($rangeAddr)
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = $fmt
}

Set-TextValue 'D2' '59.148.12'
Set-TextValue 'E2' '  +0.33%  '
Set-TextValue 'D3' '2.576.07'
Set-TextValue 'E3' '  -1.14%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '570.76'
Set-TextValue 'E5' '  +2.58%  '
Set-TextValue 'D6' '143.34'
Set-TextValue 'E6' '  -0.42%  '
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '0.595'
Set-TextValue 'E8' '  -0.46%  '
Set-TextValue 'D9' '2.579.84'
Set-TextValue 'E9' '  -1.42%  '
Set-TextValue 'D10' '6.71'
Set-TextValue 'E10' '  -1.82%  '
Set-TextValue 'E11' '  +2.86%  '
Set-TextValue 'D13' '0.346'
Set-TextValue 'E13' '  +3.01%  '
Set-TextValue 'D14' '3.024.23'
Set-TextValue 'E14' '  -1.32%  '
Set-TextValue 'D15' '59.164.91'
Set-TextValue 'E15' '  +0.42%  '
Set-TextValue 'D16' '22.46'
Set-TextValue 'E16' '  +7.45%  '
Set-TextValue 'D18' '2.584.79'
Set-TextValue 'E18' '  -1.16%  '
Set-TextValue 'E19' '  +1.45%  '
Set-TextValue 'D20' '337.98'
Set-TextValue 'D21' '10.25'
Set-TextValue 'E21' '  +1.37%  '
Set-TextValue 'E23' '  +0.10%  '
Set-TextValue 'D24' '64.51'
Set-TextValue 'E24' '  -3.26%  '
Set-TextValue 'D25' '0.458'
Set-TextValue 'E25' '  +6.79%  '
Set-TextValue 'B26' 'Binance-PegBSC-USD'
Set-TextValue 'C26' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  +0.36%  '
Set-TextValue 'B27' 'Kaspa'
Set-TextValue 'C27' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D27' '0.161'
Set-TextValue 'E27' '  -0.22%  '
Set-TextValue 'D28' '7.25'
Set-TextValue 'E28' '  +0.76%  '
Set-TextValue 'D29' '0.0₃0782'
Set-TextValue 'E29' '  +2.79%  '
Set-TextValue 'E30' '  +0.06%  '
Set-TextValue 'E31' '  +0.21%  '
Set-TextValue 'D32' '6.06'
Set-TextValue 'E32' '  +0.90%  '
Set-TextValue 'D33' '158.41'
Set-TextValue 'E33' '  +2.69%  '
Set-TextValue 'E35' '  +1.72%  '
Set-TextValue 'E36' '  +1.83%  '
Set-TextValue 'E37' '  -3.75%  '
Set-TextValue 'E38' '  -2.21%  '
Set-TextValue 'D39' '37.17'
Set-TextValue 'E39' '  +0.22%  '
Set-TextValue 'E40' '  +2.41%  '
Set-TextValue 'E41' '  +1.96%  '
Set-TextValue 'D42' '293.49'
Set-TextValue 'E42' '  +3.51%  '
Set-TextValue 'E43' '  +0.13%  '
Set-TextValue 'E44' '  +2.37%  '
Set-TextValue 'D45' '127.88'
Set-TextValue 'E45' '  +7.24%  '
Set-TextValue 'D46' '0.593'
Set-TextValue 'E46' '  -1.29%  '
Set-TextValue 'E47' '  -0.34%  '
Set-TextValue 'E48' '  +2.03%  '
Set-TextValue 'E49' '  +0.20%  '
Set-TextValue 'D50' '0.0233'
Set-TextValue 'E50' '  +1.98%  '
Set-TextValue 'D51' '1.947.56'
Set-TextValue 'E51' '  -0.42%  '
